# The species records in rows 16, 17 and 18 were shuffled: the data that
# used to live on row 17 now belongs on row 16, row 18's data moved up to
# row 17, and row 16's original data moved down to row 18. Only a subset
# of columns actually differ between the three rows (A, B, E, F, G, H, Q,
# R, Z, AB) - the rest (C, D, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW,
# AX, AY) are identical across the three rows already, so nothing else
# needs to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")
$rows = @(16, 17, 18)

# Snapshot the current values for the affected cells before overwriting
# anything. (Value2 is used for reading since it reliably returns the
# underlying scalar; Value is used below for writing.)
$values = @{}
foreach ($r in $rows) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $values[$addr] = $ws.Range($addr).Value2
    }
}

# New row 16 gets old row 17's values, new row 17 gets old row 18's
# values, new row 18 gets old row 16's values (a cyclic rotation).
$mapping = @{ 16 = 17; 17 = 18; 18 = 16 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $values["$col$src"]
    }
}
